$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM values update the edge-weight columns (G,H) and the
# receptor/edge specificity metrics (K..T) for every existing row, and the
# "Target cluster" (column D) values get reshuffled / a new cluster
# ("Inflammatory-Mac") is introduced, adding a 6th data row.

# Row 2 (D2 = ECs) numeric values
$ws.Range("G2").Value = 0.2284785
$ws.Range("H2").Value = 0.456957
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.363908
$ws.Range("N2").Value = 1.091724
$ws.Range("O2").Value = 0.01118972054940699
$ws.Range("P2").Value = 0.01663265480083782
$ws.Range("Q2").Value = 0.08314515397799999
$ws.Range("R2").Value = 0.498870923868
$ws.Range("S2").Value = 0.01118972054940699
$ws.Range("T2").Value = 0.01663265480083782

# Row 3 (D3 = FAPs) numeric values
$ws.Range("G3").Value = 0.2284785
$ws.Range("H3").Value = 0.456957
$ws.Range("N3").Value = 0.460698
$ws.Range("O3").Value = 0.004721964413781051
$ws.Range("P3").Value = 0.007018835164781924
$ws.Range("Q3").Value = 0.035086529331
$ws.Range("R3").Value = 0.210519175986
$ws.Range("S3").Value = 0.004721964413781051
$ws.Range("T3").Value = 0.007018835164781924

# Row 4: D4 changes from "MuSCs" to "Inflammatory-Mac"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.2284785
$ws.Range("H4").Value = 0.456957
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06252866666666666
$ws.Range("N4").Value = 0.187586
$ws.Range("O4").Value = 0.001922679101110775
$ws.Range("P4").Value = 0.002857913889838424
$ws.Range("Q4").Value = 0.014286455967
$ws.Range("R4").Value = 0.085718735802
$ws.Range("S4").Value = 0.001922679101110775
$ws.Range("T4").Value = 0.002857913889838424

# Row 5: D5 changes from "Resolving-Mac" to "MuSCs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.2284785
$ws.Range("H5").Value = 0.456957
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 31.927516
$ws.Range("N5").Value = 63.85503199999999
$ws.Range("O5").Value = 0.9817315966582778
$ws.Range("P5").Value = 0.9728454303033116
$ws.Range("Q5").Value = 7.294750964405999
$ws.Range("R5").Value = 29.179003857624
$ws.Range("S5").Value = 0.9817315966582778
$ws.Range("T5").Value = 0.9728454303033116

# Add new row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf5"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.2284785
$ws.Range("H6").Value = 0.456957
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01411566666666667
$ws.Range("N6").Value = 0.042347
$ws.Range("O6").Value = 0.0004340392774233579
$ws.Range("P6").Value = 0.0006451658412300904
$ws.Range("Q6").Value = 0.0032251263465
$ws.Range("R6").Value = 0.019350758079
$ws.Range("S6").Value = 0.0004340392774233579
$ws.Range("T6").Value = 0.0006451658412300904
